# Share Product 15 Test Cases
# Update the mobile number test value on the "Input" and "Output" sheets,
# clear the special highlighting style on those cells, and move the
# active-cell selection as recorded in the workbook views.

$wb = $excel.ActiveWorkbook

# --- Input sheet -----------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("B7").Value = 9901020304
$wsInput.Range("B7").Style = "Normal"
$wsInput.Range("E12").Select()

# --- Output sheet (kept active/selected last, matching the saved file) ---
$wsOutput = $wb.Worksheets.Item("Output")
$wsOutput.Range("B2").Value = 9901020304
$wsOutput.Range("B2").Style = "Normal"
$wsOutput.Range("C19").Select()
